{"js": "// Fill the header row of the first table with column labels:\n// numero | nombre | Carnet | telefono (first 4 cells of row 0).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst headers = [\"numero\", \"nombre\", \"Carnet\", \"telefono\"];\nfor (let col = 0; col < headers.length; col++) {\n  table.getCell(0, col).value = headers[col];\n}\n\nawait context.sync();\n", "ps1": "# Fill the header row of the first table with column labels:\n# numero | nombre | Carnet | telefono (first 4 cells of row 1).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$headers = @(\"numero\", \"nombre\", \"Carnet\", \"telefono\")\nfor ($c = 1; $c -le $headers.Length; $c++) {\n    $t.Cell(1, $c).Range.Text = $headers[$c - 1]\n}\n"}
